$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2625.3125
$ws.Range("I86").Value = 2621.7856
$ws.Range("J86").Value = 2650
$ws.Range("K86").Value = 2621.7856
$ws.Range("L86").Value = 2650
$ws.Range("M86").Value = -1498.7856
$ws.Range("N86").Value = -4896
$ws.Range("H89").Value = 2625.3125
$ws.Range("I89").Value = 2621.7856
$ws.Range("J89").Value = 2650
$ws.Range("K89").Value = 13108.928
$ws.Range("L89").Value = 13250
$ws.Range("M89").Value = -7492.928
$ws.Range("N89").Value = -24482
$ws.Range("I113").Value = 1532.5
$ws.Range("J113").Value = 1771.4286
$ws.Range("K113").Value = 1532.5
$ws.Range("L113").Value = 1771.4286
$ws.Range("M113").Value = 1721.5
$ws.Range("N113").Value = -8279.428599999999
$ws.Range("H116").Value = 2630.1667
$ws.Range("I116").Value = 1901.25
$ws.Range("J116").Value = 2895.2273
$ws.Range("K116").Value = 1901.25
$ws.Range("L116").Value = 2895.2273
$ws.Range("M116").Value = 1540.75
$ws.Range("N116").Value = -9779.2273
$ws.Range("H129").Value = 1788.0834
$ws.Range("I129").Value = 996.1667
$ws.Range("J129").Value = 2580
$ws.Range("K129").Value = 2988.5001
$ws.Range("L129").Value = 7740
$ws.Range("M129").Value = 2011.4999
$ws.Range("N129").Value = -17740

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5559.875
$ws.Range("I63").Value = 2160
$ws.Range("J63").Value = 7599.8
$ws.Range("K63").Value = 2160
$ws.Range("L63").Value = 7599.8
$ws.Range("M63").Value = -1474
$ws.Range("N63").Value = -8971.799999999999
$ws.Range("H66").Value = 5559.875
$ws.Range("I66").Value = 2160
$ws.Range("J66").Value = 7599.8
$ws.Range("K66").Value = 10800
$ws.Range("L66").Value = 37999
$ws.Range("M66").Value = -7368
$ws.Range("N66").Value = -44863
$ws.Range("H97").Value = 515.44446
$ws.Range("I97").Value = 515.44446
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 515.44446
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -19.44446000000005
$ws.Range("N97").Value = ""
$ws.Range("H102").Value = 1512.9375
$ws.Range("I102").Value = 1552.8572
$ws.Range("K102").Value = 1552.8572
$ws.Range("M102").Value = 69.14280000000008
$ws.Range("H122").Value = 2213.1875
$ws.Range("I122").Value = 2194.6667
$ws.Range("K122").Value = 6584.000100000001
$ws.Range("M122").Value = -4134.000100000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2440.5938
$ws.Range("I86").Value = 2345.6956
$ws.Range("J86").Value = 2683.111
$ws.Range("K86").Value = 2345.6956
$ws.Range("L86").Value = 2683.111
$ws.Range("M86").Value = -1222.6956
$ws.Range("N86").Value = -4929.111
$ws.Range("H89").Value = 2440.5938
$ws.Range("I89").Value = 2345.6956
$ws.Range("J89").Value = 2683.111
$ws.Range("K89").Value = 11728.478
$ws.Range("L89").Value = 13415.555
$ws.Range("M89").Value = -6112.477999999999
$ws.Range("N89").Value = -24647.555
$ws.Range("H94").Value = 891.2381
$ws.Range("I94").Value = 889.2353000000001
$ws.Range("J94").Value = 899.75
$ws.Range("K94").Value = 889.2353000000001
$ws.Range("L94").Value = 899.75
$ws.Range("M94").Value = -438.2353000000001
$ws.Range("N94").Value = -1801.75
$ws.Range("H99").Value = 916.8946999999999
$ws.Range("I99").Value = 936.4286
$ws.Range("K99").Value = 936.4286
$ws.Range("M99").Value = 561.5714
$ws.Range("H134").Value = 70306.875
$ws.Range("I134").Value = 86077.69500000001
$ws.Range("J134").Value = 1966.6666
$ws.Range("K134").Value = 258233.085
$ws.Range("L134").Value = 5899.9998
$ws.Range("M134").Value = -255698.085
$ws.Range("N134").Value = -10969.9998

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3608.0908
$ws.Range("I62").Value = 3455.7144
$ws.Range("J62").Value = 3874.75
$ws.Range("K62").Value = 3455.7144
$ws.Range("L62").Value = 3874.75
$ws.Range("M62").Value = -2831.7144
$ws.Range("N62").Value = -5122.75
$ws.Range("H65").Value = 3608.0908
$ws.Range("I65").Value = 3455.7144
$ws.Range("J65").Value = 3874.75
$ws.Range("K65").Value = 17278.572
$ws.Range("L65").Value = 19373.75
$ws.Range("M65").Value = -14158.572
$ws.Range("N65").Value = -25613.75
$ws.Range("H122").Value = 3379618.8
$ws.Range("I122").Value = 4808867.5
$ws.Range("K122").Value = 14426602.5
$ws.Range("M122").Value = -14424152.5
$ws.Range("H132").Value = 2240.513
$ws.Range("I132").Value = 2244.923
$ws.Range("J132").Value = 2231.6924
$ws.Range("K132").Value = 6734.768999999999
$ws.Range("L132").Value = 6695.0772
$ws.Range("M132").Value = -4204.768999999999
$ws.Range("N132").Value = -11755.0772

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 166668080
$ws.Range("I75").Value = 1831
$ws.Range("J75").Value = 333334340
$ws.Range("K75").Value = 5493
$ws.Range("L75").Value = 1000003020
$ws.Range("M75").Value = -4495
$ws.Range("N75").Value = -1000005016
$ws.Range("H78").Value = 166668080
$ws.Range("I78").Value = 1831
$ws.Range("J78").Value = 333334340
$ws.Range("K78").Value = 16479
$ws.Range("L78").Value = 3000009060
$ws.Range("M78").Value = -11487
$ws.Range("N78").Value = -3000019044

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 110581.45
$ws.Range("I80").Value = 1998.3334
$ws.Range("J80").Value = 240881.2
$ws.Range("K80").Value = 1998.3334
$ws.Range("L80").Value = 240881.2
$ws.Range("M80").Value = -1000.3334
$ws.Range("N80").Value = -242877.2
$ws.Range("H83").Value = 110581.45
$ws.Range("I83").Value = 1998.3334
$ws.Range("J83").Value = 240881.2
$ws.Range("K83").Value = 9991.666999999999
$ws.Range("L83").Value = 1204406
$ws.Range("M83").Value = -4999.666999999999
$ws.Range("N83").Value = -1214390
$ws.Range("H107").Value = 539.13043
$ws.Range("I107").Value = 331.16666
$ws.Range("J107").Value = 766
$ws.Range("K107").Value = 331.16666
$ws.Range("L107").Value = 766
$ws.Range("M107").Value = 1588.83334
$ws.Range("N107").Value = -4606
$ws.Range("H126").Value = 2705.75
$ws.Range("I126").Value = 3359.8
$ws.Range("J126").Value = 2238.5715
$ws.Range("K126").Value = 10079.4
$ws.Range("L126").Value = 6715.7145
$ws.Range("M126").Value = -7609.400000000001
$ws.Range("N126").Value = -11655.7145
$ws.Range("H132").Value = 3148.65
$ws.Range("I132").Value = 2787.25
$ws.Range("J132").Value = 3690.75
$ws.Range("K132").Value = 8361.75
$ws.Range("L132").Value = 11072.25
$ws.Range("M132").Value = -5831.75
$ws.Range("N132").Value = -16132.25

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7351.864
$ws.Range("I122").Value = 9242.666999999999
$ws.Range("J122").Value = 3300.1428
$ws.Range("K122").Value = 27728.001
$ws.Range("L122").Value = 9900.428400000001
$ws.Range("M122").Value = -25278.001
$ws.Range("N122").Value = -14800.4284

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = ""
$ws.Range("H14").Value = 9888.888999999999
$ws.Range("I14").Value = 9000
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 9000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = -8832
$ws.Range("N14").Value = -10336
$ws.Range("H62").Value = 11990.182
$ws.Range("I62").Value = 8648.833000000001
$ws.Range("J62").Value = 15999.8
$ws.Range("K62").Value = 8648.833000000001
$ws.Range("L62").Value = 15999.8
$ws.Range("M62").Value = -8024.833000000001
$ws.Range("N62").Value = -17247.8
$ws.Range("H65").Value = 11990.182
$ws.Range("I65").Value = 8648.833000000001
$ws.Range("J65").Value = 15999.8
$ws.Range("K65").Value = 43244.165
$ws.Range("L65").Value = 79999
$ws.Range("M65").Value = -40124.165
$ws.Range("N65").Value = -86239
$ws.Range("H75").Value = 22086.666
$ws.Range("J75").Value = 22086.666
$ws.Range("L75").Value = 22086.666
$ws.Range("N75").Value = -23958.666
$ws.Range("H78").Value = 22086.666
$ws.Range("J78").Value = 22086.666
$ws.Range("L78").Value = 66259.99800000001
$ws.Range("N78").Value = -75619.99800000001
$ws.Range("H96").Value = 3934.7
$ws.Range("I96").Value = 750.5
$ws.Range("J96").Value = 8711
$ws.Range("K96").Value = 750.5
$ws.Range("L96").Value = 8711
$ws.Range("M96").Value = 622.5
$ws.Range("N96").Value = -11457
